# Update the "time_taken" column (F) on the "data" sheet with refreshed query timestamps.
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$timestamps = @(
    "2021-10-05 14:20:22.409019",
    "2021-10-05 14:20:22.409026",
    "2021-10-05 14:20:22.409029",
    "2021-10-05 14:20:22.409032",
    "2021-10-05 14:20:22.409035",
    "2021-10-05 14:20:22.409037",
    "2021-10-05 14:20:22.409040",
    "2021-10-05 14:20:22.409042",
    "2021-10-05 14:20:22.409045",
    "2021-10-05 14:20:22.409048",
    "2021-10-05 14:20:22.409050",
    "2021-10-05 14:20:22.409052",
    "2021-10-05 14:20:22.409055",
    "2021-10-05 14:20:22.409057",
    "2021-10-05 14:20:22.409060",
    "2021-10-05 14:20:22.409062",
    "2021-10-05 14:20:22.409065",
    "2021-10-05 14:20:22.409068",
    "2021-10-05 14:20:22.409070",
    "2021-10-05 14:20:22.409072",
    "2021-10-05 14:20:22.409075"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $data.Range("F$row").Value = $timestamps[$i]
}

# Add a new "metadata" tab, positioned right after the "data" tab, describing
# the panel data source that was queried.
$metadata = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $data)
$metadata.Name = "metadata"

$metadata.Range("B1").Value = "data_name"
$metadata.Range("C1").Value = "data_id"
$metadata.Range("D1").Value = "data_version"
$metadata.Range("E1").Value = "data_version_created"
$metadata.Range("F1").Value = "panel_query_time"
$metadata.Range("G1").Value = "panel_get_request"

$metadata.Range("A2").Value = 0
$metadata.Range("B2").Value = "Familial Tumours Syndromes of the central & peripheral Nervous system"
$metadata.Range("C2").Value = 167
# "1.10" must stay text (not collapse to the number 1.1) and must not end up
# with a quote-prefix style, so round it through a text formula + paste-values.
$metadata.Range("D2").Formula = '="1.10"'
$metadata.Range("D2").Copy()
$metadata.Range("D2").PasteSpecial(-4163)
$metadata.Range("E2").Value = "2021-08-02T08:59:31.190391Z"
$metadata.Range("F2").Value = "2021-10-05 14:20:22.405380"
$metadata.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/167/?format=json"

# Re-apply the bold/bordered header style (matches the "data" sheet's header look)
# to the header row, and the plain numeric style used for the "data" sheet's
# leading index column to A2 - without creating any new style table entries.
$data.Range("B1").Copy()
$metadata.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$metadata.Range("A2").PasteSpecial(-4122)

[void]$metadata.Range("A1").Select()

Write-Host "metadata tab added"
